# Auto-generated Excel COM-interop script applying the Ravana_Profits leve-profit refresh.
# For each sheet, sets the updated market-price/profit cell values captured by the diff.
# Cells whose new value is $null are cleared (matching cell removal in the source diff).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (46 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 70000
$ws.Range("J81").Value = 70000
$ws.Range("L81").Value = 70000
$ws.Range("N81").Value = -71996
$ws.Range("H84").Value = 70000
$ws.Range("J84").Value = 70000
$ws.Range("L84").Value = 210000
$ws.Range("N84").Value = -219984
$ws.Range("H98").Value = 566
$ws.Range("J98").Value = 444
$ws.Range("L98").Value = 444
$ws.Range("N98").Value = -3440
$ws.Range("H112").Value = 1390.9615
$ws.Range("J112").Value = 1813.2941
$ws.Range("L112").Value = 5439.8823
$ws.Range("N112").Value = -7655.8823
$ws.Range("H121").Value = 1812.5883
$ws.Range("J121").Value = 1919.9375
$ws.Range("L121").Value = 5759.8125
$ws.Range("N121").Value = -9253.8125
$ws.Range("H122").Value = 566
$ws.Range("J122").Value = 444
$ws.Range("L122").Value = 1332
$ws.Range("N122").Value = -6232
$ws.Range("H132").Value = 796.2917
$ws.Range("I132").Value = 885.4286
$ws.Range("J132").Value = 172.33333
$ws.Range("K132").Value = 2656.2858
$ws.Range("L132").Value = 516.99999
$ws.Range("M132").Value = -126.2857999999997
$ws.Range("N132").Value = -5576.99999
$ws.Range("H135").Value = 922.6
$ws.Range("I135").Value = 821.1539
$ws.Range("K135").Value = 7390.3851
$ws.Range("M135").Value = -4855.3851
$ws.Range("H137").Value = 3224.9666
$ws.Range("I137").Value = 1764.2222
$ws.Range("K137").Value = 5292.6666
$ws.Range("M137").Value = -2742.6666
$ws.Range("H138").Value = 5329.294
$ws.Range("I138").Value = 1601.8462
$ws.Range("J138").Value = 17443.5
$ws.Range("K138").Value = 4805.5386
$ws.Range("L138").Value = 52330.5
$ws.Range("M138").Value = 334.4614000000001
$ws.Range("N138").Value = -62610.5

# ---- Sheet: ARM (28 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1022.5
$ws.Range("I2").Value = 696.6667
$ws.Range("K2").Value = 696.6667
$ws.Range("M2").Value = -583.6667
$ws.Range("H32").Value = 3205.6052
$ws.Range("I32").Value = 2580.8823
$ws.Range("K32").Value = 2580.8823
$ws.Range("M32").Value = -2293.8823
$ws.Range("H74").Value = 2459.5557
$ws.Range("I74").Value = 2535.1765
$ws.Range("K74").Value = 2535.1765
$ws.Range("M74").Value = -1661.1765
$ws.Range("H77").Value = 2459.5557
$ws.Range("I77").Value = 2535.1765
$ws.Range("K77").Value = 12675.8825
$ws.Range("M77").Value = -8307.8825
$ws.Range("H80").Value = 78110
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = $null
$ws.Range("H83").Value = 78110
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = $null
$ws.Range("H116").Value = 1022.5
$ws.Range("I116").Value = 696.6667
$ws.Range("K116").Value = 696.6667
$ws.Range("M116").Value = 1597.3333

# ---- Sheet: BSM (39 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1022.5
$ws.Range("I3").Value = 696.6667
$ws.Range("K3").Value = 696.6667
$ws.Range("M3").Value = -582.6667
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 2000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 2000
$ws.Range("M64").Value = $null
$ws.Range("N64").Value = -2450
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 2000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 2000
$ws.Range("M67").Value = $null
$ws.Range("N67").Value = -3560
$ws.Range("H99").Value = 1259.5
$ws.Range("I99").Value = 1153.7142
$ws.Range("K99").Value = 1153.7142
$ws.Range("M99").Value = 344.2858000000001
$ws.Range("H105").Value = 3643.5
$ws.Range("I105").Value = 3795.6667
$ws.Range("J105").Value = 3187
$ws.Range("K105").Value = 3795.6667
$ws.Range("L105").Value = 3187
$ws.Range("M105").Value = -2048.6667
$ws.Range("N105").Value = -6681
$ws.Range("H132").Value = 114499.5
$ws.Range("J132").Value = 114499.5
$ws.Range("L132").Value = 114499.5
$ws.Range("N132").Value = -124619.5
$ws.Range("H134").Value = 4412.3335
$ws.Range("I134").Value = 4212.125
$ws.Range("K134").Value = 12636.375
$ws.Range("M134").Value = -10101.375
$ws.Range("H138").Value = 124899.336
$ws.Range("J138").Value = 124899.336
$ws.Range("L138").Value = 124899.336
$ws.Range("N138").Value = -135179.336

# ---- Sheet: CRP (34 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2925.842
$ws.Range("I31").Value = 1895.909
$ws.Range("J31").Value = 4342
$ws.Range("K31").Value = 1895.909
$ws.Range("L31").Value = 4342
$ws.Range("M31").Value = -1600.909
$ws.Range("N31").Value = -4932
$ws.Range("H34").Value = 2925.842
$ws.Range("I34").Value = 1895.909
$ws.Range("J34").Value = 4342
$ws.Range("K34").Value = 1895.909
$ws.Range("L34").Value = 4342
$ws.Range("M34").Value = -1693.909
$ws.Range("N34").Value = -4746
$ws.Range("H74").Value = 34987.25
$ws.Range("J74").Value = 34987.25
$ws.Range("L74").Value = 34987.25
$ws.Range("N74").Value = -36735.25
$ws.Range("H77").Value = 34987.25
$ws.Range("J77").Value = 34987.25
$ws.Range("L77").Value = 104961.75
$ws.Range("N77").Value = -113697.75
$ws.Range("H105").Value = 4999.5
$ws.Range("I105").Value = 4999
$ws.Range("K105").Value = 4999
$ws.Range("M105").Value = -3252
$ws.Range("H132").Value = 3909.875
$ws.Range("I132").Value = 3754.238
$ws.Range("K132").Value = 11262.714
$ws.Range("M132").Value = -8732.714
$ws.Range("H134").Value = 3868.9443
$ws.Range("I134").Value = 3868.9443
$ws.Range("K134").Value = 11606.8329
$ws.Range("M134").Value = -9071.832900000001

# ---- Sheet: CUL (27 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4712
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = $null
$ws.Range("H64").Value = 100
$ws.Range("J64").Value = 100
$ws.Range("L64").Value = 300
$ws.Range("N64").Value = -840
$ws.Range("H66").Value = 4712
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = $null
$ws.Range("H67").Value = 100
$ws.Range("J67").Value = 100
$ws.Range("L67").Value = 300
$ws.Range("N67").Value = -2172
$ws.Range("H113").Value = 1222.1428
$ws.Range("J113").Value = 1008.8461
$ws.Range("L113").Value = 3026.5383
$ws.Range("N113").Value = -7366.5383
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = $null
$ws.Range("N117").Value = $null

# ---- Sheet: GSM (18 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7034.6665
$ws.Range("I80").Value = 1099
$ws.Range("J80").Value = 10002.5
$ws.Range("K80").Value = 1099
$ws.Range("L80").Value = 10002.5
$ws.Range("M80").Value = -101
$ws.Range("N80").Value = -11998.5
$ws.Range("H83").Value = 7034.6665
$ws.Range("I83").Value = 1099
$ws.Range("J83").Value = 10002.5
$ws.Range("K83").Value = 5495
$ws.Range("L83").Value = 50012.5
$ws.Range("M83").Value = -503
$ws.Range("N83").Value = -59996.5
$ws.Range("H132").Value = 2321.4285
$ws.Range("I132").Value = 2039.579
$ws.Range("K132").Value = 6118.737
$ws.Range("M132").Value = -3588.737

# ---- Sheet: LTW (12 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 35000
$ws.Range("I36").Value = 35000
$ws.Range("K36").Value = 35000
$ws.Range("M36").Value = -34438
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = $null
$ws.Range("H132").Value = 2889.4375
$ws.Range("I132").Value = 2173.7
$ws.Range("K132").Value = 6521.099999999999
$ws.Range("M132").Value = -3991.099999999999

# ---- Sheet: WVR (12 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 785
$ws.Range("J107").Value = 523.3333
$ws.Range("L107").Value = 1569.9999
$ws.Range("N107").Value = -5409.9999
$ws.Range("H122").Value = 1458.8
$ws.Range("J122").Value = 777
$ws.Range("L122").Value = 2331
$ws.Range("N122").Value = -7231
$ws.Range("H132").Value = 2576.6
$ws.Range("I132").Value = 1688.4286
$ws.Range("K132").Value = 5065.2858
$ws.Range("M132").Value = -2535.2858
